$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be silently parsed as a number by Excel
# (losing formatting, e.g. trailing zeros). Force them to Text first.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range('D2').Value = '43.633.48'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '2.276.90'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '122.91'
$ws.Range('E5').Value = '  +6.26%  '
$ws.Range('D6').Value = '266.34'
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('E7').Value = '  +2.47%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').Value = '48.07'
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('D11').Value = '0.0948'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '9.23'
$ws.Range('E12').Value = '  +3.29%  '
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').Value = '15.47'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '0.912'
$ws.Range('E15').Value = '  +3.23%  '
$ws.Range('D16').Value = '2.618.01'
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('D17').Value = '2.271.71'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '43.587.14'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').Value = '72.28'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').Value = '235.39'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').Value = '9.58'
$ws.Range('E24').Value = '  -3.47%  '
$ws.Range('D25').Value = '2.88'
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('E26').Value = '  +2.56%  '
$ws.Range('E27').Value = '  +1.68%  '
$ws.Range('D28').Value = '42.15'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').Value = '172.30'
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('D32').Value = '21.73'
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('D33').Value = '0.0919'
$ws.Range('E33').Value = '  -1.42%  '
$ws.Range('D34').Value = '5.74'
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '4.33'
$ws.Range('E35').Value = '  +13.70%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = '0.130'
$ws.Range('E36').Value = '  +2.00%  '
$ws.Range('D37').Value = '0.0376'
$ws.Range('E37').Value = '  +4.89%  '
$ws.Range('D38').Value = '4.62'
$ws.Range('E38').Value = '  -1.27%  '
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('E40').Value = '  +4.59%  '
$ws.Range('D41').Value = '13.96'
$ws.Range('E41').Value = '  -4.46%  '
$ws.Range('D42').Value = '74.14'
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('D43').Value = '0.239'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('D46').Value = '5.73'
$ws.Range('E46').Value = '  -10.38%  '
$ws.Range('D47').Value = '74.09'
$ws.Range('E47').Value = '  +40.76%  '
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.100'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = '8.53'
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('D51').Value = '101.70'
$ws.Range('E51').Value = '  -1.15%  '

# Restore default (unstyled) cell style now that the text is safely stored
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
